# Update countries & provincias Spain
#
# 1) Reorder country names "Hungria" / "Tailandia" (Hungria now listed
#    before Tailandia) and "Namibia" / "San Vicente y las Granadinas"
#    (Namibia now listed before San Vicente y las Granadinas).
# 2) Refresh COVID numbers (Casos totales/Nuevos casos/Casos activos/
#    Recuperados/Casos criticos/Muertes hoy/Muertes) for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the two country-name pairs -----------------------------------
# Row 63 was "Tailandia", row 64 was "Hungria" -> swap so Hungria comes first.
$ws.Range("A63").Value = "Hungria"
$ws.Range("A64").Value = "Tailandia"

# Row 194 was "San Vicente y las Granadinas", row 195 was "Namibia"
# -> swap so Namibia comes first.
$ws.Range("A194").Value = "Namibia"
$ws.Range("A195").Value = "San Vicente y las Granadinas"

# --- 2) Update the statistics -----------------------------------------------
# Row 9 - Alemania
$ws.Range("D9").Value = 130600
$ws.Range("E9").Value = 27555

# Row 24 - Mexico
$ws.Range("D24").Value = 13447
$ws.Range("E24").Value = 6580

# Row 50 - Australia
$ws.Range("E50").Value = 890
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 95

# Row 63 - now Hungria (new data for this row)
$ws.Range("B63").Value = 2998
$ws.Range("C63").Value = 56
$ws.Range("D63").Value = 629
$ws.Range("E63").Value = 2029
$ws.Range("F63").Value = 51
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 340

# Row 64 - now Tailandia (new data for this row)
$ws.Range("B64").Value = 2969
$ws.Range("C64").Value = 3
$ws.Range("D64").Value = 2739
$ws.Range("E64").Value = 176
$ws.Range("F64").Value = 61
$ws.Range("H64").Value = 54

# Row 72 - Uzbekistan
$ws.Range("B72").Value = 2127
$ws.Range("C72").Value = 9
$ws.Range("E72").Value = 847
